$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("TestData")
$ws.Activate()

# The "Result" column (S) previously cached "PASS" for the rows that had
# already been run; reset it back to blank for every row so the next
# data-driven run starts clean (rows 4 & 7 were already blank).
$ws.Range("S2").ClearContents()
$ws.Range("S3").ClearContents()
$ws.Range("S5").ClearContents()
$ws.Range("S6").ClearContents()
$ws.Range("S8").ClearContents()

# Restore the view: scrolled so column E is left-most on screen, with the
# active selection on Q11.
$excel.ActiveWindow.ScrollColumn = 5
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Q11").Select()
